$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Periodo Mora" period labels (column E, rows 16-22) with the
# new/updated periods, entered in reverse chronological order (2408 down to
# 2402) as the author did when refreshing the account statement data.
$ws.Range("E16").Value = "2408"
$ws.Range("E17").Value = "2407"
$ws.Range("E18").Value = "2406"
$ws.Range("E19").Value = "2405"
$ws.Range("E20").Value = "2404"
$ws.Range("E21").Value = "2403"
$ws.Range("E22").Value = "2402"

# Corresponding "Valor Mora" amounts follow the period values (the smaller
# 36400 debt moved from period 2408 to period 2402, swapping with 52000).
$ws.Range("F16").Value = 36400
$ws.Range("F17").Value = 52000
$ws.Range("F18").Value = 52000
$ws.Range("F19").Value = 52000
$ws.Range("F20").Value = 52000
$ws.Range("F21").Value = 52000
$ws.Range("F22").Value = 52000
